$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected. Temporarily unlock just the cells we need to
# write (per-cell, so we don't disturb the style of neighboring cells),
# write the new value, then re-lock the cell so the effective protection
# state of the sheet is unchanged.
function Set-ProtectedCellValue {
    param($range, $value)
    $range.Locked = $false
    $range.Value = $value
    $range.Locked = $true
}

# Update the confidential disclaimer text: date changes from 2021-03-25 to 2021-03-26
Set-ProtectedCellValue $ws.Range("A10") "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-26 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-7
Set-ProtectedCellValue $ws.Range("D2") 0.4913051408920505
Set-ProtectedCellValue $ws.Range("E2") 0.003946329913180779

Set-ProtectedCellValue $ws.Range("D3") 0.3316457666055673
Set-ProtectedCellValue $ws.Range("E3") 0.01829999999999998

Set-ProtectedCellValue $ws.Range("D4") 0.09328940302486449
Set-ProtectedCellValue $ws.Range("E4") 0.01722817764165385

Set-ProtectedCellValue $ws.Range("D5") 0.05527004256545857
Set-ProtectedCellValue $ws.Range("E5") 0.00300023078698386

Set-ProtectedCellValue $ws.Range("D6") 0.02848964691205902
Set-ProtectedCellValue $ws.Range("E6") 0.03425117528542643

Set-ProtectedCellValue $ws.Range("D7") 0.9999999999999999
Set-ProtectedCellValue $ws.Range("E7") 0.01075680288378744
